$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = '29.442.59'
$ws.Cells.Item(2, 5).Value = '  +0.30%  '

$ws.Cells.Item(3, 4).Value = '1.870.30'
$ws.Cells.Item(3, 5).Value = '  -0.58%  '

$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Cells.Item(5, 4) '243.69'
$ws.Cells.Item(5, 5).Value = '  +0.35%  '

$ws.Cells.Item(6, 2).Value = 'XRP'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Cells.Item(6, 4) '0.7072'
$ws.Cells.Item(6, 5).Value = '  -0.61%  '

Set-TextValue $ws.Cells.Item(7, 4) '1.001'
$ws.Cells.Item(7, 5).Value = '  -0.05%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.3157'
$ws.Cells.Item(8, 5).Value = '  +0.52%  '

Set-TextValue $ws.Cells.Item(9, 4) '0.07881'
$ws.Cells.Item(9, 5).Value = '  -1.85%  '

Set-TextValue $ws.Cells.Item(10, 4) '24.69'
$ws.Cells.Item(10, 5).Value = '  -1.69%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.07999'
$ws.Cells.Item(11, 5).Value = '  -3.93%  '

$ws.Cells.Item(12, 4).Value = '1.890.20'
$ws.Cells.Item(12, 5).Value = '  +1.09%  '

Set-TextValue $ws.Cells.Item(13, 4) '5.219'
$ws.Cells.Item(13, 5).Value = '  -0.80%  '

Set-TextValue $ws.Cells.Item(14, 4) '94.09'
$ws.Cells.Item(14, 5).Value = '  -0.73%  '

Set-TextValue $ws.Cells.Item(15, 4) '0.7055'
$ws.Cells.Item(15, 5).Value = '  -1.87%  '

Set-TextValue $ws.Cells.Item(16, 4) '6.503'
$ws.Cells.Item(16, 5).Value = '  +2.33%  '

$ws.Cells.Item(17, 4).Value = '29.497.62'
$ws.Cells.Item(17, 5).Value = '  +0.42%  '

Set-TextValue $ws.Cells.Item(18, 4) '0.000008374'
$ws.Cells.Item(18, 5).Value = '  -3.51%  '

Set-TextValue $ws.Cells.Item(19, 4) '257.71'
$ws.Cells.Item(19, 5).Value = '  +5.89%  '

$ws.Cells.Item(20, 4).Value = '2.136.07'
$ws.Cells.Item(20, 5).Value = '  -0.55%  '

Set-TextValue $ws.Cells.Item(21, 4) '13.21'
$ws.Cells.Item(21, 5).Value = '  -1.06%  '

$ws.Cells.Item(22, 5).Value = '  -0.10%  '

Set-TextValue $ws.Cells.Item(23, 4) '7.644'
$ws.Cells.Item(23, 5).Value = '  -2.56%  '

$ws.Cells.Item(24, 5).Value = '  -0.11%  '

Set-TextValue $ws.Cells.Item(25, 4) '0.1558'
$ws.Cells.Item(25, 5).Value = '  -1.14%  '

Set-TextValue $ws.Cells.Item(26, 4) '9.076'
$ws.Cells.Item(26, 5).Value = '  -0.19%  '

Set-TextValue $ws.Cells.Item(27, 4) '161.07'
$ws.Cells.Item(27, 5).Value = '  -1.50%  '

Set-TextValue $ws.Cells.Item(28, 4) '18.84'
$ws.Cells.Item(28, 5).Value = '  +1.20%  '

$ws.Cells.Item(29, 5).Value = '  -0.65%  '

Set-TextValue $ws.Cells.Item(30, 4) '4.338'
$ws.Cells.Item(30, 5).Value = '  -2.27%  '

$ws.Cells.Item(31, 5).Value = '  -2.12%  '

Set-TextValue $ws.Cells.Item(32, 4) '1.207'
$ws.Cells.Item(32, 5).Value = '  +0.12%  '

Set-TextValue $ws.Cells.Item(33, 4) '0.05325'

Set-TextValue $ws.Cells.Item(34, 4) '1.900'
$ws.Cells.Item(34, 5).Value = '  -2.17%  '

$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Cells.Item(35, 4) '1.175'
$ws.Cells.Item(35, 5).Value = '  -0.36%  '

$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(36, 4) '0.7492'
$ws.Cells.Item(36, 5).Value = '  -4.02%  '

Set-TextValue $ws.Cells.Item(37, 4) '2.712'

Set-TextValue $ws.Cells.Item(38, 4) '0.01883'
$ws.Cells.Item(38, 5).Value = '  -0.29%  '

$ws.Cells.Item(39, 4).Value = '1.265.15'
$ws.Cells.Item(39, 5).Value = '  -0.73%  '

$ws.Cells.Item(40, 5).Value = '  +0.24%  '

Set-TextValue $ws.Cells.Item(41, 4) '0.8999'
$ws.Cells.Item(41, 5).Value = '  -2.29%  '

Set-TextValue $ws.Cells.Item(42, 4) '108.68'
$ws.Cells.Item(42, 5).Value = '  -4.41%  '

$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(43, 4) '71.93'
$ws.Cells.Item(43, 5).Value = '  -3.47%  '

$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Cells.Item(44, 4) '5.976'
$ws.Cells.Item(44, 5).Value = '  -8.68%  '

$ws.Cells.Item(45, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(45, 4) '0.00000000130'
$ws.Cells.Item(45, 5).Value = '  +2.08%  '

$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Cells.Item(46, 4) '1.000'
$ws.Cells.Item(46, 5).Value = '  -0.08%  '

$ws.Cells.Item(47, 4).Value = '2.034.42'
$ws.Cells.Item(47, 5).Value = '  -0.26%  '

Set-TextValue $ws.Cells.Item(48, 4) '1.798'
$ws.Cells.Item(48, 5).Value = '  -0.68%  '

Set-TextValue $ws.Cells.Item(49, 4) '0.5196'
$ws.Cells.Item(49, 5).Value = '  -0.49%  '

Set-TextValue $ws.Cells.Item(50, 4) '9.517'
$ws.Cells.Item(50, 5).Value = '  -0.45%  '

Set-TextValue $ws.Cells.Item(51, 4) '0.4333'
